# Fixed data, fixed error
# - Rename the "IsAvailable" column header to "Availability" (sheet CurrentRelay, D1)
# - Fix the D2 data value: it was incorrectly stored as boolean TRUE, should be numeric 0
# - Update the active selection on the CurrentRelay sheet to J7 (was P26)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CurrentRelay")

# Rename header D1 from "IsAvailable" to "Availability"
$ws.Range("D1").Value = "Availability"

# Fix D2: was boolean TRUE, should be numeric 0
$ws.Range("D2").Value = 0

# Update selection to reflect the author's current cursor position
[void]$ws.Range("J7").Select()
